$d = $word.ActiveDocument

# --- Change 1 & 2: mark the two inline picture runs as "not spell-checked"
# (Word's re-insertion of a pasted/auto-inserted picture run sets <w:noProof/>
# on the run's rPr) ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# --- Change 3: append a new paragraph (several runs) after the last paragraph,
# right before the sectPr ---
$insertPoint = $d.Content.End
$insRange = $d.Range($insertPoint, $insertPoint)
$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">Sørg for at have referencebillede når det er fuldstændigt tømt for </w:t></w:r><w:r><w:t xml:space="preserve">tracer i 2D projektionerne!!! Så du kan bruge det samme på tværs. Sørg også for at </w:t></w:r><w:r><w:t xml:space="preserve">skylle helt igennem med vand. Tager måske lang tid? Min high pressure video giver ingen mening. Den der ligger derinde nu og hedder </w:t></w:r><w:r><w:t>non-average_THISONE</w:t></w:r><w:r><w:t>, det er high pressure, hvor jeg har brugt low pressur</w:t></w:r><w:r><w:t xml:space="preserve">res </w:t></w:r><w:r><w:t>reference</w:t></w:r><w:r><w:t>-billede som reference til high pressure også. De mindste værdier er negativ</w:t></w:r><w:r><w:t>. Jeg har ikke absolut-værdi. Og det er img-ref, dvs. img er større end ref ved de første frames. Det er jo vildt – dvs. img har allerede haft mere tracer ved frame 1</w:t></w:r><w:r><w:t xml:space="preserve"> end ref havde. Den er IKKE TØMT for tracer!! Så vi får ikke den første front med. Og faktisk så er de sidste frames højere. Dvs. </w:t></w:r><w:r><w:t>img er højere end ref. Dvs. det er endnu mere tømt det sidste… HMMM. Interessant.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$insRange.InsertXML($newParaXml)
